$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the two brand-new rows (114, 115) the same date-formatted style as column A
# elsewhere in the sheet (copy format only from A113, which already carries it).
$ws.Range("A113").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)

# Rewrite rows 93-115 (dates shift by one day starting 2021-02-08, and a new
# row is appended at the end for 2021-02-21). Values are written directly since
# this sheet stores only static numbers (no formulas).
$ws.Cells.Item(93,1).Value = 44235
$ws.Cells.Item(93,2).Value = 5
$ws.Cells.Item(93,3).Value = 23
$ws.Cells.Item(93,4).Value = 151.7550804961731
$ws.Cells.Item(94,1).Value = 44236
$ws.Cells.Item(94,2).Value = 9
$ws.Cells.Item(94,3).Value = 28
$ws.Cells.Item(94,4).Value = 184.7453153866456
$ws.Cells.Item(95,1).Value = 44237
$ws.Cells.Item(95,2).Value = 1
$ws.Cells.Item(95,3).Value = 29
$ws.Cells.Item(95,4).Value = 191.34336236474
$ws.Cells.Item(96,1).Value = 44238
$ws.Cells.Item(96,2).Value = 2
$ws.Cells.Item(96,3).Value = 31
$ws.Cells.Item(96,4).Value = 204.539456320929
$ws.Cells.Item(97,1).Value = 44239
$ws.Cells.Item(97,2).Value = 7
$ws.Cells.Item(97,3).Value = 32
$ws.Cells.Item(97,4).Value = 211.1375032990235
$ws.Cells.Item(98,1).Value = 44240
$ws.Cells.Item(98,2).Value = 3
$ws.Cells.Item(98,3).Value = 32
$ws.Cells.Item(98,4).Value = 211.1375032990235
$ws.Cells.Item(99,1).Value = 44241
$ws.Cells.Item(99,2).Value = 4
$ws.Cells.Item(99,3).Value = 31
$ws.Cells.Item(99,4).Value = 204.539456320929
$ws.Cells.Item(100,1).Value = 44242
$ws.Cells.Item(100,2).Value = 6
$ws.Cells.Item(100,3).Value = 31
$ws.Cells.Item(100,4).Value = 204.539456320929
$ws.Cells.Item(101,1).Value = 44243
$ws.Cells.Item(101,2).Value = 9
$ws.Cells.Item(101,3).Value = 38
$ws.Cells.Item(101,4).Value = 250.7257851675904
$ws.Cells.Item(102,1).Value = 44244
$ws.Cells.Item(102,2).Value = 0
$ws.Cells.Item(102,3).Value = 40
$ws.Cells.Item(102,4).Value = 263.9218791237794
$ws.Cells.Item(103,1).Value = 44245
$ws.Cells.Item(103,2).Value = 2
$ws.Cells.Item(103,3).Value = 39
$ws.Cells.Item(103,4).Value = 257.3238321456849
$ws.Cells.Item(104,1).Value = 44246
$ws.Cells.Item(104,2).Value = 14
$ws.Cells.Item(104,3).Value = 35
$ws.Cells.Item(104,4).Value = 230.9316442333069
$ws.Cells.Item(105,1).Value = 44247
$ws.Cells.Item(105,2).Value = 5
$ws.Cells.Item(105,3).Value = 36
$ws.Cells.Item(105,4).Value = 237.5296912114014
$ws.Cells.Item(106,1).Value = 44248
$ws.Cells.Item(106,2).Value = 3
$ws.Cells.Item(106,3).Value = 36
$ws.Cells.Item(106,4).Value = 237.5296912114014
$ws.Cells.Item(107,1).Value = 44249
$ws.Cells.Item(107,2).Value = 2
$ws.Cells.Item(107,3).Value = 34
$ws.Cells.Item(107,4).Value = 224.3335972552125
$ws.Cells.Item(108,1).Value = 44250
$ws.Cells.Item(108,2).Value = 10
$ws.Cells.Item(108,3).Value = 34
$ws.Cells.Item(108,4).Value = 224.3335972552125
$ws.Cells.Item(109,1).Value = 44251
$ws.Cells.Item(109,2).Value = 0
$ws.Cells.Item(109,3).Value = 34
$ws.Cells.Item(109,4).Value = 224.3335972552125
$ws.Cells.Item(110,1).Value = 44252
$ws.Cells.Item(110,2).Value = 0
$ws.Cells.Item(110,3).Value = 32
$ws.Cells.Item(110,4).Value = 211.1375032990235
$ws.Cells.Item(111,1).Value = 44253
$ws.Cells.Item(111,2).Value = 14
$ws.Cells.Item(111,3).Value = 34
$ws.Cells.Item(111,4).Value = 224.3335972552125
$ws.Cells.Item(112,1).Value = 44254
$ws.Cells.Item(112,2).Value = 5
$ws.Cells.Item(112,3).Value = 31
$ws.Cells.Item(112,4).Value = 204.539456320929
$ws.Cells.Item(113,1).Value = 44255
$ws.Cells.Item(113,2).Value = 1
$ws.Cells.Item(114,1).Value = 44256
$ws.Cells.Item(114,2).Value = 4
$ws.Cells.Item(115,1).Value = 44257
$ws.Cells.Item(115,2).Value = 7

# The 7-day rolling-sum columns (C, D) also change for a few rows just before
# the inserted date, because the new data point enters their trailing window.
$ws.Cells.Item(90,3).Value = 22
$ws.Cells.Item(90,4).Value = 145.1570335180786
$ws.Cells.Item(91,3).Value = 28
$ws.Cells.Item(91,4).Value = 184.7453153866456
$ws.Cells.Item(92,3).Value = 27
$ws.Cells.Item(92,4).Value = 178.1472684085511
